# Updated remaining queries for C3DC
# Replace "id"/"study.id"/"participant.id" style join keys with the fully
# qualified "study_id"/"participant_id" keys across every SQL query stored
# on Sheet1, and widen column C to fit the (now slightly longer) text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Update-Query([string]$text) {
    $text = $text.Replace('std.id = prt."study.id"', 'std.study_id = prt."study.study_id"')
    $text = $text.Replace('prt.id = dgn."participant.id"', 'prt.participant_id = dgn."participant.participant_id"')
    $text = $text.Replace('prt.id = trt."participant.id"', 'prt.participant_id = trt."participant.participant_id"')
    $text = $text.Replace('prt.id = trr."participant.id"', 'prt.participant_id = trr."participant.participant_id"')
    $text = $text.Replace('prt.id = srv."participant.id"', 'prt.participant_id = srv."participant.participant_id"')
    $text = $text.Replace('std.id = rfs."study.id"', 'std.study_id = rfs."study.study_id"')
    return $text
}

# Cells holding a SQL query that needs the join-key update.
$cells = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")

foreach ($addr in $cells) {
    $rng = $ws.Range($addr)
    $rng.Value2 = Update-Query $rng.Value2
}

# Column C width: was bestFit-computed 60.83203125, now an explicit 67.5
# (no more auto-fit). Excel's ColumnWidth property is offset from the raw
# OOXML "width" (character units) by the default ~5px margin (5/6 of a
# character at the workbook's standard font), so back that out here to
# land on exactly 67.5 in the saved file.
$ws.Columns.Item(3).ColumnWidth = 67.5 - (5/6)
